$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = "Where the advertisement being applied to be built will be located"
$ws.Range("B27").Value = "How long the proposed advertisement will be shown."
$ws.Range("B29").Value = "What type of advertisements are proposed and how many there will be."
$ws.Range("B33").Value = "Name and contact information if an agent is being used."
$ws.Range("B37").Value = "Name and contact information if an agent is being used."
$ws.Range("B45").Value = "Telephone number and email address of the applicant."
$ws.Range("B49").Value = "Name and contact information for the parties making the application."
$ws.Range("B55").Value = "Checking whether all the requirements of the form have been met, such as proof of payment or supporting documentation."
$ws.Range("B56").Value = "What community consultation activities have taken place as part of the application"
$ws.Range("B58").Value = "Details of any conflict of interest that may exist between the applicant and planning authority."
$ws.Range("B61").Value = "Names and contact details for all parties with an interest in the proposed develpoment."
$ws.Range("B79").Value = "Signed and dated verification of the application's accuracy."
$ws.Range("B82").Value = "Details of pre-application advice received from the planning authority"
$ws.Range("B87").Value = "Details of the proposed advertisements such as their size and how they are made"
$ws.Range("B99").Value = "Where the proposed development will be built."
$ws.Range("B108").Value = "Information to help the planning authority arrange a site visit"
